# Solved Problems: Group Anagrams, Insert Intervals and Merge Intervals.
#
# Adds four new "Array" topic rows (45-48: First Missing Positive, Group
# Anagrams, Merged Intervals, Insert Intervals) plus two trailing note-only
# rows (49-50: Gas Station, Candy) to the tracker, and moves the sheet's
# viewport/selection down to the newly added area.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the new rows by pasting row 44's formatting down into
# rows 45:48 - this reuses the existing cell styles (center-aligned text,
# date-formatted column F) instead of minting brand-new style entries.
$ws.Range("A44:I44").Copy()
$ws.Range("A45:I48").PasteSpecial(-4122)

# Row 45 - First Missing Positive
$ws.Cells.Item(45, 1).Value = 44
$ws.Cells.Item(45, 2).Value = "Array "
$ws.Cells.Item(45, 3).Value = "First Missing Positive"
$ws.Cells.Item(45, 4).Value = "Hard"
$ws.Cells.Item(45, 5).Value = "Done"
$ws.Cells.Item(45, 6).Value = 45887
$ws.Cells.Item(45, 7).Value = "O(n)"
$ws.Cells.Item(45, 8).Value = "O(1)"
$ws.Cells.Item(45, 9).Value = "Cycle Sort"

# Row 46 - Group Anagrams
$ws.Cells.Item(46, 1).Value = 45
$ws.Cells.Item(46, 2).Value = "Array"
$ws.Cells.Item(46, 3).Value = "Group Anagrams"
$ws.Cells.Item(46, 4).Value = "Medium"
$ws.Cells.Item(46, 5).Value = "Done"
$ws.Cells.Item(46, 6).Value = 45887
$ws.Cells.Item(46, 7).Value = "O(NK)"
$ws.Cells.Item(46, 8).Value = "O(NK)"
$ws.Cells.Item(46, 9).Value = "Frequency Count + Map"

# Row 47 - Merged Intervals
$ws.Cells.Item(47, 1).Value = 46
$ws.Cells.Item(47, 2).Value = "Array"
$ws.Cells.Item(47, 3).Value = "Merged Intervals"
$ws.Cells.Item(47, 4).Value = "Medium"
$ws.Cells.Item(47, 5).Value = "Done"
$ws.Cells.Item(47, 6).Value = 45889
$ws.Cells.Item(47, 7).Value = "O(n * log n)"
$ws.Cells.Item(47, 8).Value = "O(n)"
$ws.Cells.Item(47, 9).Value = "Sorting"

# Row 48 - Insert Intervals
$ws.Cells.Item(48, 1).Value = 47
$ws.Cells.Item(48, 2).Value = "Array"
$ws.Cells.Item(48, 3).Value = "Insert Intervals"
$ws.Cells.Item(48, 4).Value = "Medium"
$ws.Cells.Item(48, 5).Value = "Done"
$ws.Cells.Item(48, 6).Value = 45889
$ws.Cells.Item(48, 7).Value = "O(n)"
$ws.Cells.Item(48, 8).Value = "O(n)"
$ws.Cells.Item(48, 9).Value = "Normal Traversing"

# Rows 49-50 - trailing notes (Column C only)
$ws.Cells.Item(49, 3).Value = "Gas Station"
$ws.Cells.Item(50, 3).Value = "Candy"

# Scroll the viewport to the newly added rows and move the selection.
$excel.ActiveWindow.ScrollRow = 32
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C51").Select()
